$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column F (the empty "intepretive themes" column), shifting
# column G ("description") left into its place.
$ws.Columns("F").Delete()

# Update the selection to match the saved cursor position in the file.
$ws.Range("J4").Select()
